$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1 (paragraph "به طور کلی Property به دو دسته تقسیم می‌شوند:"):
# insert the word "ها " right before "به دو دسته", splitting the old single
# trailing run into three runs: " " | "ها " | "به دو دسته تقسیم می‌شوند:"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(9)
$r1 = $p1.Range.Duplicate
$r1.Find.Execute("به دو دسته", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPos1 = $r1.Start
$rins1 = $d.Range($insPos1, $insPos1)
$rins1.Text = "ها "
# Force the freshly inserted "ها " text to become its own run (distinct from
# the run that follows it) by toggling a character property on/off.
$rNew1 = $d.Range($insPos1, $insPos1 + 3)
$rNew1.Bold = $true
$rNew1.Bold = $false

# ---------------------------------------------------------------------------
# Change 2 (paragraph "Property های متناظر فیلد (که باید هم get و هم set‌
# داشته باشند)"): insert the word "یک " right before "فیلد", splitting the
# run into three runs: " های متناظر " | "یک " | "فیلد (که باید هم "
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(10)
$r2 = $p2.Range.Duplicate
$r2.Find.Execute("فیلد (که باید هم", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPos2 = $r2.Start
$rins2 = $d.Range($insPos2, $insPos2)
$rins2.Text = "یک "
$rNew2 = $d.Range($insPos2, $insPos2 + 3)
$rNew2.Bold = $true
$rNew2.Bold = $false

# ---------------------------------------------------------------------------
# Change 3 (same paragraph): "باشند)" -> "باشد)" (drop the "ن"), splitting
# the final run into two runs: "‌ داشته باش" | "د)"
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(10)
$r3 = $p3.Range.Duplicate
$r3.Find.Execute("باشند)", $true, $false, $false, $false, $false, $true, 1, $false, "باشد)", 2)
$rSplit3 = $d.Range($r3.Start + 3, $r3.End)
$rSplit3.Bold = $true
$rSplit3.Bold = $false

# ---------------------------------------------------------------------------
# Change 4/5 (paragraph "Property های محاسباتی که معمولا فقط get دارند."):
# the run before the _GoBack bookmark gains "تی که معمولا" and the run after
# it loses the same text, i.e. only the <w:t> contents of the two existing
# runs change -- the bookmark itself stays put between them.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(11)
$r4a = $p4.Range.Duplicate
$r4a.Find.Execute(" های محاسبا", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4a.Text = " های محاسباتی که معمولا"

$p4b = $d.Paragraphs(11)
$r4b = $p4b.Range.Duplicate
$r4b.Find.Execute("تی که معمولا فقط ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4b.Text = " فقط "
